$wb = $excel.ActiveWorkbook

# --- Rename headers on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after "Monthly Trend" (last sheet) ---
# Copy "Monthly Trend" so the new sheet inherits the same sheetPr/pageMargins
# layout, then wipe it clean before writing the forecast data.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsMonthly.Copy($null, $lastSheet)
$wsForecast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast.Name = "PO Forecast"
$wsForecast.Cells.Clear()

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy header style (bold + border) from Weekly Quantity header row
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Copy date-number-format style from Weekly Quantity date column
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A52").PasteSpecial(-4122)

# Data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$wsForecast.Range("A2").Value = 45088.99999999999
$wsForecast.Range("B2").Value = 90
$wsForecast.Range("C2").Value = -33.25236437160159
$wsForecast.Range("D2").Value = 220.6674205371612
$wsForecast.Range("A3").Value = 45095.99999999999
$wsForecast.Range("B3").Value = 91
$wsForecast.Range("C3").Value = -31.91542602258824
$wsForecast.Range("D3").Value = 224.851713709209
$wsForecast.Range("A4").Value = 45102.99999999999
$wsForecast.Range("B4").Value = 91
$wsForecast.Range("C4").Value = -41.23173078064999
$wsForecast.Range("D4").Value = 227.1747460595184
$wsForecast.Range("A5").Value = 45109.99999999999
$wsForecast.Range("B5").Value = 91
$wsForecast.Range("C5").Value = -38.69296146977539
$wsForecast.Range("D5").Value = 220.4060665962255
$wsForecast.Range("A6").Value = 45116.99999999999
$wsForecast.Range("B6").Value = 91
$wsForecast.Range("C6").Value = -33.73672087699821
$wsForecast.Range("D6").Value = 228.6584032668584
$wsForecast.Range("A7").Value = 45123.99999999999
$wsForecast.Range("B7").Value = 92
$wsForecast.Range("C7").Value = -41.43448430340104
$wsForecast.Range("D7").Value = 211.965493688066
$wsForecast.Range("A8").Value = 45130.99999999999
$wsForecast.Range("B8").Value = 92
$wsForecast.Range("C8").Value = -40.94629211573053
$wsForecast.Range("D8").Value = 227.3502351537119
$wsForecast.Range("A9").Value = 45137.99999999999
$wsForecast.Range("B9").Value = 92
$wsForecast.Range("C9").Value = -48.18715614308186
$wsForecast.Range("D9").Value = 220.3784015785853
$wsForecast.Range("A10").Value = 45144.99999999999
$wsForecast.Range("B10").Value = 92
$wsForecast.Range("C10").Value = -40.77672182856819
$wsForecast.Range("D10").Value = 217.2838788568291
$wsForecast.Range("A11").Value = 45151.99999999999
$wsForecast.Range("B11").Value = 93
$wsForecast.Range("C11").Value = -31.3001521305664
$wsForecast.Range("D11").Value = 220.5003057701362
$wsForecast.Range("A12").Value = 45172.99999999999
$wsForecast.Range("B12").Value = 93
$wsForecast.Range("C12").Value = -33.54482060380943
$wsForecast.Range("D12").Value = 232.1520979593188
$wsForecast.Range("A13").Value = 45186.99999999999
$wsForecast.Range("B13").Value = 94
$wsForecast.Range("C13").Value = -35.49635396274866
$wsForecast.Range("D13").Value = 226.7906797415848
$wsForecast.Range("A14").Value = 45193.99999999999
$wsForecast.Range("B14").Value = 94
$wsForecast.Range("C14").Value = -33.09130004531143
$wsForecast.Range("D14").Value = 230.462493828867
$wsForecast.Range("A15").Value = 45200.99999999999
$wsForecast.Range("B15").Value = 94
$wsForecast.Range("C15").Value = -41.97529343938159
$wsForecast.Range("D15").Value = 228.789476228436
$wsForecast.Range("A16").Value = 45207.99999999999
$wsForecast.Range("B16").Value = 95
$wsForecast.Range("C16").Value = -19.12144564813279
$wsForecast.Range("D16").Value = 219.6460973409561
$wsForecast.Range("A17").Value = 45214.99999999999
$wsForecast.Range("B17").Value = 95
$wsForecast.Range("C17").Value = -32.05681715675873
$wsForecast.Range("D17").Value = 227.2685538729481
$wsForecast.Range("A18").Value = 45221.99999999999
$wsForecast.Range("B18").Value = 95
$wsForecast.Range("C18").Value = -40.1877547167898
$wsForecast.Range("D18").Value = 223.9900154976339
$wsForecast.Range("A19").Value = 45228.99999999999
$wsForecast.Range("B19").Value = 95
$wsForecast.Range("C19").Value = -41.93215800842254
$wsForecast.Range("D19").Value = 225.1770286217988
$wsForecast.Range("A20").Value = 45235.99999999999
$wsForecast.Range("B20").Value = 96
$wsForecast.Range("C20").Value = -37.39979019063497
$wsForecast.Range("D20").Value = 228.6907893206156
$wsForecast.Range("A21").Value = 45242.99999999999
$wsForecast.Range("B21").Value = 96
$wsForecast.Range("C21").Value = -28.50755862837952
$wsForecast.Range("D21").Value = 227.030766376339
$wsForecast.Range("A22").Value = 45256.99999999999
$wsForecast.Range("B22").Value = 96
$wsForecast.Range("C22").Value = -42.04238470360409
$wsForecast.Range("D22").Value = 228.329829917114
$wsForecast.Range("A23").Value = 45270.99999999999
$wsForecast.Range("B23").Value = 97
$wsForecast.Range("C23").Value = -40.11891778454307
$wsForecast.Range("D23").Value = 218.8061385654881
$wsForecast.Range("A24").Value = 45277.99999999999
$wsForecast.Range("B24").Value = 97
$wsForecast.Range("C24").Value = -36.59276096823627
$wsForecast.Range("D24").Value = 226.3554438222206
$wsForecast.Range("A25").Value = 45298.99999999999
$wsForecast.Range("B25").Value = 98
$wsForecast.Range("C25").Value = -30.69977444147855
$wsForecast.Range("D25").Value = 224.4875342451169
$wsForecast.Range("A26").Value = 45312.99999999999
$wsForecast.Range("B26").Value = 98
$wsForecast.Range("C26").Value = -30.62612354942763
$wsForecast.Range("D26").Value = 231.2479668084241
$wsForecast.Range("A27").Value = 45326.99999999999
$wsForecast.Range("B27").Value = 99
$wsForecast.Range("C27").Value = -17.12660937909204
$wsForecast.Range("D27").Value = 228.4065192944457
$wsForecast.Range("A28").Value = 45333.99999999999
$wsForecast.Range("B28").Value = 99
$wsForecast.Range("C28").Value = -27.71627632977506
$wsForecast.Range("D28").Value = 225.0584451772117
$wsForecast.Range("A29").Value = 45347.99999999999
$wsForecast.Range("B29").Value = 100
$wsForecast.Range("C29").Value = -31.58394465952206
$wsForecast.Range("D29").Value = 228.1712903215081
$wsForecast.Range("A30").Value = 45361.99999999999
$wsForecast.Range("B30").Value = 100
$wsForecast.Range("C30").Value = -27.85371120020351
$wsForecast.Range("D30").Value = 225.8514423541627
$wsForecast.Range("A31").Value = 45417.99999999999
$wsForecast.Range("B31").Value = 102
$wsForecast.Range("C31").Value = -24.65261950648859
$wsForecast.Range("D31").Value = 224.8592971201275
$wsForecast.Range("A32").Value = 45424.99999999999
$wsForecast.Range("B32").Value = 102
$wsForecast.Range("C32").Value = -30.58760715960308
$wsForecast.Range("D32").Value = 234.9429293626072
$wsForecast.Range("A33").Value = 45459.99999999999
$wsForecast.Range("B33").Value = 104
$wsForecast.Range("C33").Value = -32.97994195952269
$wsForecast.Range("D33").Value = 239.9687004273711
$wsForecast.Range("A34").Value = 45466.99999999999
$wsForecast.Range("B34").Value = 104
$wsForecast.Range("C34").Value = -35.88211247138101
$wsForecast.Range("D34").Value = 228.3675892781152
$wsForecast.Range("A35").Value = 45473.99999999999
$wsForecast.Range("B35").Value = 104
$wsForecast.Range("C35").Value = -27.39212352922847
$wsForecast.Range("D35").Value = 233.2041589894545
$wsForecast.Range("A36").Value = 45557.99999999999
$wsForecast.Range("B36").Value = 107
$wsForecast.Range("C36").Value = -23.28292227926661
$wsForecast.Range("D36").Value = 237.8478368335042
$wsForecast.Range("A37").Value = 45564.99999999999
$wsForecast.Range("B37").Value = 107
$wsForecast.Range("C37").Value = -28.60244897631751
$wsForecast.Range("D37").Value = 242.5454662211633
$wsForecast.Range("A38").Value = 45571.99999999999
$wsForecast.Range("B38").Value = 108
$wsForecast.Range("C38").Value = -13.72399553664427
$wsForecast.Range("D38").Value = 239.5645710008765
$wsForecast.Range("A39").Value = 45578.99999999999
$wsForecast.Range("B39").Value = 108
$wsForecast.Range("C39").Value = -17.65896452028817
$wsForecast.Range("D39").Value = 235.7520335980122
$wsForecast.Range("A40").Value = 45585.99999999999
$wsForecast.Range("B40").Value = 108
$wsForecast.Range("C40").Value = -16.19829368883916
$wsForecast.Range("D40").Value = 240.1687150419598
$wsForecast.Range("A41").Value = 45592.99999999999
$wsForecast.Range("B41").Value = 108
$wsForecast.Range("C41").Value = -23.63241671296524
$wsForecast.Range("D41").Value = 244.9615502843103
$wsForecast.Range("A42").Value = 45599.99999999999
$wsForecast.Range("B42").Value = 109
$wsForecast.Range("C42").Value = -19.48190100453304
$wsForecast.Range("D42").Value = 237.3348562231876
$wsForecast.Range("A43").Value = 45613.99999999999
$wsForecast.Range("B43").Value = 109
$wsForecast.Range("C43").Value = -17.06923634362372
$wsForecast.Range("D43").Value = 245.295209211424
$wsForecast.Range("A44").Value = 45641.99999999999
$wsForecast.Range("B44").Value = 110
$wsForecast.Range("C44").Value = -17.00037661326446
$wsForecast.Range("D44").Value = 243.7745284000646
$wsForecast.Range("A45").Value = 45648.99999999999
$wsForecast.Range("B45").Value = 110
$wsForecast.Range("C45").Value = -13.85064095228687
$wsForecast.Range("D45").Value = 235.3609081628308
$wsForecast.Range("A46").Value = 45655.99999999999
$wsForecast.Range("B46").Value = 111
$wsForecast.Range("C46").Value = -26.15242708620317
$wsForecast.Range("D46").Value = 232.0013830929126
$wsForecast.Range("A47").Value = 45662.99999999999
$wsForecast.Range("B47").Value = 111
$wsForecast.Range("C47").Value = -30.34275003488068
$wsForecast.Range("D47").Value = 233.6107039858079
$wsForecast.Range("A48").Value = 45669.99999999999
$wsForecast.Range("B48").Value = 111
$wsForecast.Range("C48").Value = -17.52167584291213
$wsForecast.Range("D48").Value = 234.9399837954702
$wsForecast.Range("A49").Value = 45676.99999999999
$wsForecast.Range("B49").Value = 111
$wsForecast.Range("C49").Value = -20.32755877344994
$wsForecast.Range("D49").Value = 243.8004571346848
$wsForecast.Range("A50").Value = 45683.99999999999
$wsForecast.Range("B50").Value = 112
$wsForecast.Range("C50").Value = -14.9964256775463
$wsForecast.Range("D50").Value = 240.1713469606066
$wsForecast.Range("A51").Value = 45690.99999999999
$wsForecast.Range("B51").Value = 112
$wsForecast.Range("C51").Value = -14.89058295702297
$wsForecast.Range("D51").Value = 244.8990305425641
$wsForecast.Range("A52").Value = 45697.99999999999
$wsForecast.Range("B52").Value = 112
$wsForecast.Range("C52").Value = -10.78122306767881
$wsForecast.Range("D52").Value = 237.2001061729267

Write-Host "PO Forecast sheet created with data"